# Bump the ValueSet metadata to version 1.1.0, as published on 2023-07-10.
#
# The "Metadata" worksheet is a two-column Property/Value table; this walks
# column A looking for the "Version" and "Date" property rows and updates the
# matching cell in column B, so the edit is resilient to the exact row
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$newVersion = "1.1.0"
$newDate = "2023-07-10T23:08:03+02:00"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) {
    $lastRow = 14
}

for ($row = 1; $row -le $lastRow; $row++) {
    $label = $ws.Cells.Item($row, 1).Value2

    if ($label -eq "Version") {
        $ws.Cells.Item($row, 2).Value = $newVersion
    }
    elseif ($label -eq "Date") {
        $ws.Cells.Item($row, 2).Value = $newDate
    }
}
